{"js": "// The intro paragraph originally read:\n//   \"T\u00f4i t\u00ean l\u00e0 Phan Thanh H\u1ea3i \u0111\u1ebfn t\u1eeb tr\u01b0\u1eddng \u0110H Kinh t\u1ebf TP. H\u1ed3 Ch\u00ed Minh. ...\"\n// The edit replaces \"\u0110H Kinh t\u1ebf TP. H\u1ed3 Ch\u00ed Minh\" with \"\u0110H C\u00f4ng ngh\u1ec7 th\u00f4ng\n// tin\" and relocates the \"_GoBack\" bookmark from the end of the document\n// (after the \"T\u1ed5 tr\u01b0\u1edfng/ T\u1ed5 ph\u00f3\" paragraph) to right after the newly\n// inserted university name, splitting the paragraph into three runs:\n//   1) \"T\u00f4i t\u00ean l\u00e0 Phan Thanh H\u1ea3i \u0111\u1ebfn t\u1eeb tr\u01b0\u1eddng \"\n//   2) \"\u0110H C\u00f4ng ngh\u1ec7 th\u00f4ng tin\"\n//      <bookmarkStart/bookmarkEnd name=\"_GoBack\">\n//   3) \". Hi\u1ec7n t\u1ea1i t\u00f4i \u0111ang l\u00e0m \u0111\u1ec1 t\u00e0i ...\"\n\nconst oldText = \"\u0110H Kinh t\u1ebf TP. H\u1ed3 Ch\u00ed Minh\";\nconst newText = \"\u0110H C\u00f4ng ngh\u1ec7 th\u00f4ng tin\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found: \" + oldText);\n}\n\nconst target = results.items[0];\n\n// The \"_GoBack\" bookmark previously sat after the \"T\u1ed5 tr\u01b0\u1edfng/ T\u1ed5 ph\u00f3\"\n// paragraph; remove it there first since it is moving.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Mark the boundary just before the text we are about to replace with a\n// throwaway bookmark. This forces the surrounding (identically-formatted)\n// text to stay in its own run instead of merging with the replacement text,\n// matching the run split seen in the target document.\nconst startRange = target.getRange(\"Start\");\nstartRange.insertBookmark(\"_splitmark_\");\nawait context.sync();\n\n// Replace the old university name with the new one; it now lands in its own\n// run because of the split point inserted above.\ntarget.insertText(newText, \"Replace\");\nawait context.sync();\n\n// Drop the throwaway bookmark and put \"_GoBack\" right after the newly\n// inserted text.\ncontext.document.deleteBookmark(\"_splitmark_\");\nawait context.sync();\n\nconst newResults = body.search(newText, { matchCase: true });\nnewResults.load(\"items\");\nawait context.sync();\n\nif (newResults.items.length === 0) {\n  throw new Error(\"Replacement text not found after insert: \" + newText);\n}\n\nconst afterRange = newResults.items[0].getRange(\"After\");\nafterRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The intro paragraph originally read:\n#   \"T\u00f4i t\u00ean l\u00e0 Phan Thanh H\u1ea3i \u0111\u1ebfn t\u1eeb tr\u01b0\u1eddng \u0110H Kinh t\u1ebf TP. H\u1ed3 Ch\u00ed Minh. ...\"\n# This edit replaces \"\u0110H Kinh t\u1ebf TP. H\u1ed3 Ch\u00ed Minh\" with \"\u0110H C\u00f4ng ngh\u1ec7 th\u00f4ng\n# tin\" and relocates the \"_GoBack\" bookmark from the end of the document\n# (right after the \"T\u1ed5 tr\u01b0\u1edfng/ T\u1ed5 ph\u00f3\" paragraph) to right after the newly\n# inserted university name, splitting the paragraph into three runs:\n#   1) \"T\u00f4i t\u00ean l\u00e0 Phan Thanh H\u1ea3i \u0111\u1ebfn t\u1eeb tr\u01b0\u1eddng \"\n#   2) \"\u0110H C\u00f4ng ngh\u1ec7 th\u00f4ng tin\"\n#      <bookmarkStart/bookmarkEnd name=\"_GoBack\">\n#   3) \". Hi\u1ec7n t\u1ea1i t\u00f4i \u0111ang l\u00e0m \u0111\u1ec1 t\u00e0i ...\"\n\n$d = $word.ActiveDocument\n\n$oldText = \"\u0110H Kinh t\u1ebf TP. H\u1ed3 Ch\u00ed Minh\"\n$newText = \"\u0110H C\u00f4ng ngh\u1ec7 th\u00f4ng tin\"\n\n# The \"_GoBack\" bookmark previously sat after the \"T\u1ed5 tr\u01b0\u1edfng/ T\u1ed5 ph\u00f3\"\n# paragraph; remove it there first since it is moving.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Locate the text to replace.\n$r = $d.Content\n$found = $r.Find.Execute($oldText)\nif (-not $found) {\n  throw \"Target text not found: $oldText\"\n}\n\n# Drop a throwaway bookmark collapsed at the start of the found range. This\n# forces the preceding (identically-formatted) text to stay in its own run\n# instead of merging with the replacement text, matching the run split seen\n# in the target document.\n$splitPoint = $d.Range($r.Start, $r.Start)\n$d.Bookmarks.Add(\"_splitmark_\", $splitPoint)\n\n# Replace the old university name with the new one; it lands in its own run\n# because of the split point inserted above. After the assignment, $r is\n# collapsed to cover exactly the newly inserted text.\n$r.Text = $newText\n\n# Drop the throwaway bookmark.\n$d.Bookmarks.Item(\"_splitmark_\").Delete()\n\n# Place \"_GoBack\" collapsed right after the newly inserted text.\n$endPoint = $d.Range($r.End, $r.End)\n$d.Bookmarks.Add(\"_GoBack\", $endPoint)\n"}
